$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Generation 0-4 (rows 2-6): Fitness 7293 -> 8357
$ws.Range("C2:C6").Value = 8357

# Generation 5-24 (rows 7-26): Fitness 7293 -> 7785
$ws.Range("C7:C26").Value = 7785
